$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary header values (labels unchanged, only amounts updated) ---
$ws.Range("E11").Value = 604010
$ws.Range("C13").Value = 8
$ws.Range("F13").Value = 11

# --- Make room for the new worker rows ---
# Old table: header row 15, data rows 16-21 (last row 21 has the thicker
# bottom-border style), footer rows 26-27.
# New table: header row 15, data rows 16-28 (last row 28 keeps the thicker
# bottom-border style), footer rows 33-34.
# Insert 7 blank rows right before the old closing row so the closing
# row (with its special formatting) ends up at row 28.
$ws.Rows("22:28").Insert()

# Re-apply the regular data-row formatting (borders/number formats) to the
# freshly inserted rows by copying the format of an existing normal row.
$ws.Range("B16:J16").Copy()
$ws.Range("B22:J28").PasteSpecial(-4122)

# Re-apply the special "last row" formatting (thicker bottom border) to the
# new final row of the table.
$ws.Range("B21:J21").Copy()
$ws.Range("B28:J28").PasteSpecial(-4122)

# --- Fill in the worker/period rows (B..G); columns H..J stay blank ---
$rows = @(
    @{ Row = 16; Doc = "1047447559"; Name = "DAVID SALOMON MORALES GONZALES"; Periodo = "2103"; Mora = 36341;  Salario = 908526 },
    @{ Row = 17; Doc = "1047447559"; Name = "DAVID SALOMON MORALES GONZALES"; Periodo = "2102"; Mora = 36341;  Salario = 908526 },
    @{ Row = 18; Doc = "1143386760"; Name = "ANDRES FELIPE FRIAS FIGUEROA";    Periodo = "1712"; Mora = 30000;  Salario = 750000 },
    @{ Row = 19; Doc = "1023908096"; Name = "LEYDI CONSTANZA PIRACHICAN DAZA"; Periodo = "2507"; Mora = 62000;  Salario = 1550000 },
    @{ Row = 20; Doc = "1023908096"; Name = "LEYDI CONSTANZA PIRACHICAN DAZA"; Periodo = "2506"; Mora = 62000;  Salario = 1550000 },
    @{ Row = 21; Doc = "1023908096"; Name = "LEYDI CONSTANZA PIRACHICAN DAZA"; Periodo = "2504"; Mora = 62000;  Salario = 1550000 },
    @{ Row = 22; Doc = "1023908096"; Name = "LEYDI CONSTANZA PIRACHICAN DAZA"; Periodo = "2503"; Mora = 62000;  Salario = 1550000 },
    @{ Row = 23; Doc = "1023908096"; Name = "LEYDI CONSTANZA PIRACHICAN DAZA"; Periodo = "2502"; Mora = 62000;  Salario = 1550000 },
    @{ Row = 24; Doc = "1050969971"; Name = "YESSICA PAOLA CASTILLA TORRES";   Periodo = "2007"; Mora = 18200;  Salario = 910000 },
    @{ Row = 25; Doc = "1002202656"; Name = "WILFRAN JOSE MENDOZA GOMEZ";      Periodo = "2402"; Mora = 52000;  Salario = 1300000 },
    @{ Row = 26; Doc = "1127608729"; Name = "VICTOR ANTONIO BERRIO TERAN";     Periodo = "2501"; Mora = 52800;  Salario = 1320000 },
    @{ Row = 27; Doc = "1066870078"; Name = "NICOLLE MARCELA MORENO MORENO";   Periodo = "2507"; Mora = 56940;  Salario = 1423500 },
    @{ Row = 28; Doc = "1043961211"; Name = "ALEXANDER GUZMAN MARTINEZ";       Periodo = "2506"; Mora = 11388;  Salario = 1423500 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("B$n").Value = "CC"
    $ws.Range("C$n").Value = $r.Doc
    $ws.Range("D$n").Value = $r.Name
    $ws.Range("E$n").Value = $r.Periodo
    $ws.Range("F$n").Value = $r.Mora
    $ws.Range("G$n").Value = $r.Salario
}
